$wb = $excel.ActiveWorkbook

# --- Labels sheet: add three new rows describing new placeholders \V, \F, \T ---
$labels = $wb.Worksheets.Item("Labels")

$labels.Range("A31").Value = "\\V"
$labels.Range("B31").Value = "event value"

$labels.Range("A32").Value = "\\F"
$labels.Range("B32").Value = "event value interpreted as temperature in Fahrenheit converted to the current temperature mode"

$labels.Range("A33").Value = "\\T"
$labels.Range("B33").Value = "event value interpreted as temperature in Celsius converted to the current temperature mode"

$labels.Range("A31").Font.Italic = $true
$labels.Range("A32").Font.Italic = $true
$labels.Range("A33").Font.Italic = $true

# --- Commands sheet: tweak wording of the placeholder note (add a comma) ---
$commands = $wb.Worksheets.Item("Commands")
$commands.Range("A3").Value = "tn:Note: The placeholders {ET}, {BT}, {time}, {ETB}, {BTB}, and {WEIGHTin} will be substituted by the current ET, BT, time, ET background, BT background value, and batch size (in g) in Serial/Artisan/CallProgram/MODBUS/S7/WebSocket commands\n"
